# Remove the extra blank paragraph that sits right before the final
# blank paragraph at the end of the document (both are empty paragraphs
# with the same bold/Segoe UI formatting that trail the last table).
$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($count - 1)
$following = $d.Paragraphs.Item($count)

$r = $d.Range($target.Range.Start, $following.Range.Start)
$r.Delete()
